$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.632.26"
$ws.Range("E2").Value = "  +0.00%  "
$ws.Range("D3").Value = "1.597.75"
$ws.Range("E3").Value = "  +0.37%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").Value = "'211.38"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.28%  "
$ws.Range("E6").Value = "  +0.40%  "
$ws.Range("E7").Value = "  +0.10%  "
$ws.Range("E8").Value = "  +0.13%  "
$ws.Range("E9").Value = "  +0.03%  "
$ws.Range("D10").Value = "'19.45"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.75%  "
$ws.Range("E11").Value = "  +0.34%  "
$ws.Range("D12").Value = "1.821.71"
$ws.Range("E12").Value = "  +0.32%  "
$ws.Range("D13").Value = "1.602.48"
$ws.Range("E13").Value = "  +0.70%  "
$ws.Range("E14").Value = "  -0.12%  "
$ws.Range("E15").Value = "  -0.31%  "
$ws.Range("E16").Value = "  -0.45%  "
$ws.Range("D17").Value = "26.633.09"
$ws.Range("E17").Value = "  +0.11%  "
$ws.Range("E18").Value = "  +0.83%  "
$ws.Range("E19").Value = "  +0.10%  "
$ws.Range("D20").Value = "'208.65"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.13%  "
$ws.Range("D21").Value = "'7.05"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +5.34%  "
$ws.Range("E22").Value = "  +0.07%  "
$ws.Range("E23").Value = "  -0.54%  "
$ws.Range("E24").Value = "  +0.02%  "
$ws.Range("D25").Value = "'145.54"
$ws.Range("D25").Style = "Normal"
$ws.Range("E26").Value = "  +0.02%  "
$ws.Range("D27").Value = "'7.15"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.19%  "
$ws.Range("E28").Value = "  +0.78%  "
$ws.Range("D29").Value = "'15.25"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.35%  "
$ws.Range("E30").Value = "  +0.53%  "
$ws.Range("D31").Value = "'1.16"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.01%  "
$ws.Range("E32").Value = "  -0.27%  "
$ws.Range("E33").Value = "  +0.91%  "
$ws.Range("D34").Value = "'0.624"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -6.17%  "
$ws.Range("D35").Value = "1.271.30"
$ws.Range("E35").Value = "  -1.88%  "
$ws.Range("E36").Value = "  +0.40%  "
$ws.Range("E37").Value = "  +0.07%  "
$ws.Range("D38").Value = "'0.0170"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.78%  "
$ws.Range("D39").Value = "'0.840"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.25%  "
$ws.Range("E40").Value = "  +2.65%  "
$ws.Range("E41").Value = "  +0.93%  "
$ws.Range("E42").Value = "  -0.62%  "
$ws.Range("D43").Value = "'64.19"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.68%  "
$ws.Range("D44").Value = "'0.944"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +17.45%  "
$ws.Range("D45").Value = "1.735.24"
$ws.Range("E45").Value = "  +0.37%  "
$ws.Range("D46").Value = "'89.98"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.07%  "
$ws.Range("E47").Value = "  +0.16%  "
$ws.Range("D48").Value = "'0.103"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +4.40%  "
$ws.Range("D49").Value = "'0.0509"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.19%  "
$ws.Range("D50").Value = "'7.47"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.10%  "
$ws.Range("E51").Value = "  +0.05%  "
